$d = $word.ActiveDocument

# Paragraph 1 ("Following sealed evidence was submitted along with the
# request of DPO, Jhelum for Comparison of Cartridge Cases with Submitted
# Firearm and Functionality Testing.") is now populated from the template
# when the report is created, so trim the boilerplate case-specific tail
# back off, leaving only "Following sealed evidence was".
$found = $d.Content.Find.Execute(
    " submitted along with the request of DPO, Jhelum for Comparison of Cartridge Cases with Submitted Firearm and Functionality Testing.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 1)

if (-not $found) {
    throw "Could not find the report boilerplate text to remove."
}
